# The deck ships two theme parts:
#   ppt/theme/theme1.xml  - "Integral" / "Red Violet"  (used by the slide master -> the
#                            visually-active theme for every slide)
#   ppt/theme/theme2.xml  - "Office Theme" / "Office"   (used only by the notes master)
#
# The target edit swaps the two themes' content: the slide master's theme becomes the
# plain default "Office Theme" colour set, and the notes master's theme becomes the
# "Integral"/"Red Violet" colour set that used to be active.
#
# The font scheme and format (fill/line/effect) scheme are byte-for-byte identical
# between the two themes already, so the only thing that actually needs to change is
# the 12-colour scheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink).
#
# PowerPoint exposes those 12 colours per-slide via Slide.ThemeColorScheme, which edits
# the slide master's theme part (theme1.xml) in place - this is the change that is
# visible across every slide in the deck, so we apply it here using the Design tab's
# underlying object model (ThemeColorScheme / RGBColor).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Index order matches MsoThemeColorSchemeIndex:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3, 8 accent4,
# 9 accent5, 10 accent6, 11 hlink, 12 folHlink
# Target values come from the standard Office theme palette.
$officeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

for ($i = 1; $i -le 12; $i++) {
    $tcs.Item($i).RGB = $officeColors[$i - 1]
}
